$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 'total (1)(2)'
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 4.212800105155133
$ws.Range("D4").Value = 2.635001223733333
$ws.Range("E4").Value = 1.951854323411999
$ws.Range("F4").Value = 2.838887807342168
$ws.Range("G4").Value = 4.624632123611441

$ws.Range("A5").Value = 'homens'
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 5.10572950074503
$ws.Range("D5").Value = 3.392140686652902
$ws.Range("E5").Value = 2.595524914670427
$ws.Range("F5").Value = 3.801519962445005
$ws.Range("G5").Value = 6.080495728501765

$ws.Range("A6").Value = 'mulheres'
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 5.976035449443504
$ws.Range("D6").Value = 3.717758752778127
$ws.Range("E6").Value = 2.602428979737852
$ws.Range("F6").Value = 4.062545578680676
$ws.Range("G6").Value = 6.403746462083348

$ws.Range("A7").Value = 'branca'
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 5.726065832786598
$ws.Range("D7").Value = 3.939974055301707
$ws.Range("E7").Value = 2.938339659139324
$ws.Range("F7").Value = 4.778322935662588
$ws.Range("G7").Value = 6.922629892138545

$ws.Range("A8").Value = 'preta ou parda'
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 5.705446586612587
$ws.Range("D8").Value = 3.372393094285793
$ws.Range("E8").Value = 2.475324705729576
$ws.Range("F8").Value = 3.470792534197381
$ws.Range("G8").Value = 5.787911953506636

$ws.Range("A9").Value = 'preta'
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 11.59420420385286
$ws.Range("D9").Value = 8.254530162702418
$ws.Range("E9").Value = 6.020645937615219
$ws.Range("F9").Value = 8.472887156131044
$ws.Range("G9").Value = 12.72869714918794

$ws.Range("A10").Value = 'parda'
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 6.269788609528693
$ws.Range("D10").Value = 3.611151093744206
$ws.Range("E10").Value = 2.620857062664615
$ws.Range("F10").Value = 3.70488056845653
$ws.Range("G10").Value = 6.455475564808933

$ws.Range("A11").Value = '10 a 24 anos'
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 6.512351399680265
$ws.Range("D11").Value = 3.745550271173594
$ws.Range("E11").Value = 2.969955259599499
$ws.Range("F11").Value = 4.394125690378891
$ws.Range("G11").Value = 6.744170318686604

$ws.Range("A12").Value = '25 a 34 anos'
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 8.617844902262037
$ws.Range("D12").Value = 5.037468080623166
$ws.Range("E12").Value = 3.523781822775497
$ws.Range("F12").Value = 5.727643083251569
$ws.Range("G12").Value = 9.393881301644262

$ws.Range("A13").Value = '35 a 39 anos'
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 12.91227528174228
$ws.Range("D13").Value = 8.877088725297837
$ws.Range("E13").Value = 5.760457042673352
$ws.Range("F13").Value = 9.956089404303633
$ws.Range("G13").Value = 14.2743366273082

$ws.Range("A14").Value = '40 a 49 anos'
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 10.03861598948733
$ws.Range("D14").Value = 6.934886081310327
$ws.Range("E14").Value = 4.517776043263018
$ws.Range("F14").Value = 7.645523494819397
$ws.Range("G14").Value = 13.06542392046357

$ws.Range("A15").Value = '50 anos ou mais'
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 8.446764596621506
$ws.Range("D15").Value = 6.827322627645295
$ws.Range("E15").Value = 4.636796862814598
$ws.Range("F15").Value = 6.674058136258643
$ws.Range("G15").Value = 10.23442828843285

$ws.Range("A16").Value = 'sem instrução'
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 13.75439957694691
$ws.Range("D16").Value = 13.32794913282037
$ws.Range("E16").Value = 8.568556187260329
$ws.Range("F16").Value = 9.851259494434519
$ws.Range("G16").Value = 16.73173330548401

$ws.Range("A17").Value = 'fundamental incompleto ou equivalente'
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 6.337282548281007
$ws.Range("D17").Value = 4.666277752136621
$ws.Range("E17").Value = 3.587023540864393
$ws.Range("F17").Value = 4.877425230356435
$ws.Range("G17").Value = 8.194365027673522

$ws.Range("A18").Value = 'fundamental completo ou equivalente'
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 11.74168047059585
$ws.Range("D18").Value = 7.653675734502754
$ws.Range("E18").Value = 5.22962705758344
$ws.Range("F18").Value = 7.769137384681653
$ws.Range("G18").Value = 13.44351621113692

$ws.Range("A19").Value = 'médio incompleto ou equivalente'
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 11.69317065244053
$ws.Range("D19").Value = 6.401524844937679
$ws.Range("E19").Value = 4.895246003140828
$ws.Range("F19").Value = 7.052949928265715
$ws.Range("G19").Value = 12.08475005952241

$ws.Range("A20").Value = 'médio completo ou equivalente'
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 8.018584491046074
$ws.Range("D20").Value = 4.418447300527082
$ws.Range("E20").Value = 3.23784051870862
$ws.Range("F20").Value = 5.604293695193507
$ws.Range("G20").Value = 8.336675918460278

$ws.Range("A21").Value = 'superior incompleto ou equivalente'
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 15.05908040531419
$ws.Range("D21").Value = 8.671496447839361
$ws.Range("E21").Value = 6.632862385581006
$ws.Range("F21").Value = 12.85224438905852
$ws.Range("G21").Value = 16.7921570772392

$ws.Range("A22").Value = 'superior completo ou equivalente'
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 17.19988235328317
$ws.Range("D22").Value = 9.43566541459598
$ws.Range("E22").Value = 6.871321971186088
$ws.Range("F22").Value = 13.37856941745204
$ws.Range("G22").Value = 13.41010361493089

$ws.Range("A23").Value = 'total (3)(4)'
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 4.229274194575787
$ws.Range("D23").Value = 2.641423975816012
$ws.Range("E23").Value = 1.96418273510481
$ws.Range("F23").Value = 2.847862094175835
$ws.Range("G23").Value = 4.637764564074551

$ws.Range("A24").Value = 'sem rendimento a menos de 1/4 do salário mínimo (3) (5)'
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = 11.70588994020808
$ws.Range("D24").Value = 10.46138064628093
$ws.Range("E24").Value = 5.410321809420506
$ws.Range("F24").Value = 8.318024315393188
$ws.Range("G24").Value = 15.46289827995497

$ws.Range("A25").Value = '1/4 a menos de 1/2 salário mínimo (3)'
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 9.137900997593514
$ws.Range("D25").Value = 6.378427299373929
$ws.Range("E25").Value = 4.613287576433751
$ws.Range("F25").Value = 5.64019762150286
$ws.Range("G25").Value = 10.79866603599297

$ws.Range("A26").Value = '1/2 a menos de 1 salário mínimo (3)'
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 7.281393581008858
$ws.Range("D26").Value = 4.704373816026497
$ws.Range("E26").Value = 3.734214053085317
$ws.Range("F26").Value = 4.926772344862705
$ws.Range("G26").Value = 9.089104552444068

$ws.Range("A27").Value = '1 a menos de 2 salários mínimos (3) '
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 8.09799328245062
$ws.Range("D27").Value = 4.569855531381092
$ws.Range("E27").Value = 3.77508433057174
$ws.Range("F27").Value = 5.82457818758136
$ws.Range("G27").Value = 8.413642646921948

$ws.Range("A28").Value = '2 salários mínimos ou mais (3) '
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 8.875113870938092
$ws.Range("D28").Value = 5.911385268256816
$ws.Range("E28").Value = 4.256984801968604
$ws.Range("F28").Value = 8.458586967953563
$ws.Range("G28").Value = 9.32558790765031

$ws.Range("A29:G35").Clear()

$ws.Range("A1:G28").Select()